$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 121.2130545
$ws.Range("H2").Value = 242.426109
$ws.Range("I2").Value = 0.4742742514031324
$ws.Range("J2").Value = 0.4028886461111009
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 121.2130545
$ws.Range("N2").Value = 242.426109
$ws.Range("O2").Value = 0.4742742514031324
$ws.Range("P2").Value = 0.4028886461111009
$ws.Range("Q2").Value = 14692.60458121997
$ws.Range("R2").Value = 58770.41832487988
$ws.Range("S2").Value = 0.2249360655440016
$ws.Range("T2").Value = 0.1623192611652359

$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 121.2130545
$ws.Range("H3").Value = 242.426109
$ws.Range("I3").Value = 0.4742742514031324
$ws.Range("J3").Value = 0.4028886461111009
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.9119573333333332
$ws.Range("N3").Value = 2.735872
$ws.Range("O3").Value = 0.003568245048871889
$ws.Range("P3").Value = 0.004546753526507615
$ws.Range("Q3").Value = 110.541133947008
$ws.Range("R3").Value = 663.2468036820479
$ws.Range("S3").Value = 0.001692326749376649
$ws.Range("T3").Value = 0.001831835372495526

$ws.Range("E4").Value = 2
$ws.Range("G4").Value = 121.2130545
$ws.Range("H4").Value = 242.426109
$ws.Range("I4").Value = 0.4742742514031324
$ws.Range("J4").Value = 0.4028886461111009
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.7040186666666667
$ws.Range("N4").Value = 2.112056
$ws.Range("O4").Value = 0.002754636680714656
$ws.Range("P4").Value = 0.003510031926267592
$ws.Range("Q4").Value = 85.33625301168401
$ws.Range("R4").Value = 512.017518070104
$ws.Range("S4").Value = 0.001306453249633553
$ws.Range("T4").Value = 0.00141415201058069

$ws.Range("E5").Value = 2
$ws.Range("G5").Value = 121.2130545
$ws.Range("H5").Value = 242.426109
$ws.Range("I5").Value = 0.4742742514031324
$ws.Range("J5").Value = 0.4028886461111009
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.6855039999999999
$ws.Range("N5").Value = 2.056512
$ws.Range("O5").Value = 0.002682193743693282
$ws.Range("P5").Value = 0.00341772319330189
$ws.Range("Q5").Value = 83.09203371196799
$ws.Range("R5").Value = 498.5522022718079
$ws.Range("S5").Value = 0.001272095429908296
$ws.Range("T5").Value = 0.001376961870131907

$ws.Range("E6").Value = 2
$ws.Range("G6").Value = 121.2130545
$ws.Range("H6").Value = 242.426109
$ws.Range("I6").Value = 0.4742742514031324
$ws.Range("J6").Value = 0.4028886461111009
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 88.26666266666666
$ws.Range("N6").Value = 264.799988
$ws.Range("O6").Value = 0.3453638350486923
$ws.Range("P6").Value = 0.4400718598158738
$ws.Range("Q6").Value = 10699.07179234778
$ws.Range("R6").Value = 64194.43075408669
$ws.Range("S6").Value = 0.1637971743294334
$ws.Range("T6").Value = 0.1772999557928116

$ws.Range("E7").Value = 2
$ws.Range("G7").Value = 121.2130545
$ws.Range("H7").Value = 242.426109
$ws.Range("I7").Value = 0.4742742514031324
$ws.Range("J7").Value = 0.4028886461111009
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 43.794673
$ws.Range("N7").Value = 87.58934600000001
$ws.Range("O7").Value = 0.1713568380748954
$ws.Range("P7").Value = 0.1455649854269483
$ws.Range("Q7").Value = 5308.486085158679
$ws.Range("R7").Value = 21233.94434063472
$ws.Range("S7").Value = 0.08127013610077881
$ws.Range("T7").Value = 0.05864647989984531

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.9119573333333332
$ws.Range("H8").Value = 2.735872
$ws.Range("I8").Value = 0.003568245048871889
$ws.Range("J8").Value = 0.004546753526507615
$ws.Range("K8").Value = 2
$ws.Range("M8").Value = 121.2130545
$ws.Range("N8").Value = 242.426109
$ws.Range("O8").Value = 0.4742742514031324
$ws.Range("P8").Value = 0.4028886461111009
$ws.Range("Q8").Value = 110.541133947008
$ws.Range("R8").Value = 663.2468036820479
$ws.Range("S8").Value = 0.001692326749376649
$ws.Range("T8").Value = 0.001831835372495526

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.9119573333333332
$ws.Range("H9").Value = 2.735872
$ws.Range("I9").Value = 0.003568245048871889
$ws.Range("J9").Value = 0.004546753526507615
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.9119573333333332
$ws.Range("N9").Value = 2.735872
$ws.Range("O9").Value = 0.003568245048871889
$ws.Range("P9").Value = 0.004546753526507615
$ws.Range("Q9").Value = 0.8316661778204442
$ws.Range("R9").Value = 7.484995600383998
$ws.Range("S9").Value = [double]"1.273237272879875E-05"
$ws.Range("T9").Value = [double]"2.067296763080943E-05"

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.9119573333333332
$ws.Range("H10").Value = 2.735872
$ws.Range("I10").Value = 0.003568245048871889
$ws.Range("J10").Value = 0.004546753526507615
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.7040186666666667
$ws.Range("N10").Value = 2.112056
$ws.Range("O10").Value = 0.002754636680714656
$ws.Range("P10").Value = 0.003510031926267592
$ws.Range("Q10").Value = 0.6420349858702221
$ws.Range("R10").Value = 5.778314872831999
$ws.Range("S10").Value = [double]"9.829218697400969E-06"
$ws.Range("T10").Value = [double]"1.595925003891149E-05"

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.9119573333333332
$ws.Range("H11").Value = 2.735872
$ws.Range("I11").Value = 0.003568245048871889
$ws.Range("J11").Value = 0.004546753526507615
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.6855039999999999
$ws.Range("N11").Value = 2.056512
$ws.Range("O11").Value = 0.002682193743693282
$ws.Range("P11").Value = 0.00341772319330189
$ws.Range("Q11").Value = 0.6251503998293332
$ws.Range("R11").Value = 5.626353598463998
$ws.Range("S11").Value = [double]"9.570724546048711E-06"
$ws.Range("T11").Value = [double]"1.553954498177224E-05"

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.9119573333333332
$ws.Range("H12").Value = 2.735872
$ws.Range("I12").Value = 0.003568245048871889
$ws.Range("J12").Value = 0.004546753526507615
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 88.26666266666666
$ws.Range("N12").Value = 264.799988
$ws.Range("O12").Value = 0.3453638350486923
$ws.Range("P12").Value = 0.4400718598158738
$ws.Range("Q12").Value = 80.49543030772621
$ws.Range("R12").Value = 724.4588727695359
$ws.Range("S12").Value = 0.001232342794471904
$ws.Range("T12").Value = 0.002000898280534589

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.9119573333333332
$ws.Range("H13").Value = 2.735872
$ws.Range("I13").Value = 0.003568245048871889
$ws.Range("J13").Value = 0.004546753526507615
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 43.794673
$ws.Range("N13").Value = 87.58934600000001
$ws.Range("O13").Value = 0.1713568380748954
$ws.Range("P13").Value = 0.1455649854269483
$ws.Range("Q13").Value = 39.93887320328533
$ws.Range("R13").Value = 239.633239219712
$ws.Range("S13").Value = 0.0006114431890510878
$ws.Range("T13").Value = 0.0006618481108260066

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.7040186666666667
$ws.Range("H14").Value = 2.112056
$ws.Range("I14").Value = 0.002754636680714656
$ws.Range("J14").Value = 0.003510031926267592
$ws.Range("K14").Value = 2
$ws.Range("M14").Value = 121.2130545
$ws.Range("N14").Value = 242.426109
$ws.Range("O14").Value = 0.4742742514031324
$ws.Range("P14").Value = 0.4028886461111009
$ws.Range("Q14").Value = 85.33625301168401
$ws.Range("R14").Value = 512.017518070104
$ws.Range("S14").Value = 0.001306453249633553
$ws.Range("T14").Value = 0.00141415201058069

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.7040186666666667
$ws.Range("H15").Value = 2.112056
$ws.Range("I15").Value = 0.002754636680714656
$ws.Range("J15").Value = 0.003510031926267592
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 0.9119573333333332
$ws.Range("N15").Value = 2.735872
$ws.Range("O15").Value = 0.003568245048871889
$ws.Range("P15").Value = 0.004546753526507615
$ws.Range("Q15").Value = 0.6420349858702221
$ws.Range("R15").Value = 5.778314872831999
$ws.Range("S15").Value = [double]"9.829218697400969E-06"
$ws.Range("T15").Value = [double]"1.595925003891149E-05"

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.7040186666666667
$ws.Range("H16").Value = 2.112056
$ws.Range("I16").Value = 0.002754636680714656
$ws.Range("J16").Value = 0.003510031926267592
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.7040186666666667
$ws.Range("N16").Value = 2.112056
$ws.Range("O16").Value = 0.002754636680714656
$ws.Range("P16").Value = 0.003510031926267592
$ws.Range("Q16").Value = 0.4956422830151112
$ws.Range("R16").Value = 4.460780547135999
$ws.Range("S16").Value = [double]"7.58802324273866E-06"
$ws.Range("T16").Value = [double]"1.232032412341778E-05"

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.7040186666666667
$ws.Range("H17").Value = 2.112056
$ws.Range("I17").Value = 0.002754636680714656
$ws.Range("J17").Value = 0.003510031926267592
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 0.6855039999999999
$ws.Range("N17").Value = 2.056512
$ws.Range("O17").Value = 0.002682193743693282
$ws.Range("P17").Value = 0.00341772319330189
$ws.Range("Q17").Value = 0.4826076120746666
$ws.Range("R17").Value = 4.343468508671999
$ws.Range("S17").Value = [double]"7.388469271160881E-06"
$ws.Range("T17").Value = [double]"1.199631752363486E-05"

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 0.7040186666666667
$ws.Range("H18").Value = 2.112056
$ws.Range("I18").Value = 0.002754636680714656
$ws.Range("J18").Value = 0.003510031926267592
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 88.26666266666666
$ws.Range("N18").Value = 264.799988
$ws.Range("O18").Value = 0.3453638350486923
$ws.Range("P18").Value = 0.4400718598158738
$ws.Range("Q18").Value = 62.14137816170311
$ws.Range("R18").Value = 559.272403455328
$ws.Range("S18").Value = 0.0009513518882174138
$ws.Range("T18").Value = 0.001544666277805673

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 0.7040186666666667
$ws.Range("H19").Value = 2.112056
$ws.Range("I19").Value = 0.002754636680714656
$ws.Range("J19").Value = 0.003510031926267592
$ws.Range("K19").Value = 2
$ws.Range("M19").Value = 43.794673
$ws.Range("N19").Value = 87.58934600000001
$ws.Range("O19").Value = 0.1713568380748954
$ws.Range("P19").Value = 0.1455649854269483
$ws.Range("Q19").Value = 30.83226729256267
$ws.Range("R19").Value = 184.993603755376
$ws.Range("S19").Value = 0.0004720258316523888
$ws.Range("T19").Value = 0.0005109377461952651

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 0.6855039999999999
$ws.Range("H20").Value = 2.056512
$ws.Range("I20").Value = 0.002682193743693282
$ws.Range("J20").Value = 0.00341772319330189
$ws.Range("K20").Value = 2
$ws.Range("M20").Value = 121.2130545
$ws.Range("N20").Value = 242.426109
$ws.Range("O20").Value = 0.4742742514031324
$ws.Range("P20").Value = 0.4028886461111009
$ws.Range("Q20").Value = 83.09203371196799
$ws.Range("R20").Value = 498.5522022718079
$ws.Range("S20").Value = 0.001272095429908296
$ws.Range("T20").Value = 0.001376961870131907

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 0.6855039999999999
$ws.Range("H21").Value = 2.056512
$ws.Range("I21").Value = 0.002682193743693282
$ws.Range("J21").Value = 0.00341772319330189
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 0.9119573333333332
$ws.Range("N21").Value = 2.735872
$ws.Range("O21").Value = 0.003568245048871889
$ws.Range("P21").Value = 0.004546753526507615
$ws.Range("Q21").Value = 0.6251503998293332
$ws.Range("R21").Value = 5.626353598463998
$ws.Range("S21").Value = [double]"9.570724546048711E-06"
$ws.Range("T21").Value = [double]"1.553954498177224E-05"

$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 0.6855039999999999
$ws.Range("H22").Value = 2.056512
$ws.Range("I22").Value = 0.002682193743693282
$ws.Range("J22").Value = 0.00341772319330189
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 0.7040186666666667
$ws.Range("N22").Value = 2.112056
$ws.Range("O22").Value = 0.002754636680714656
$ws.Range("P22").Value = 0.003510031926267592
$ws.Range("Q22").Value = 0.4826076120746666
$ws.Range("R22").Value = 4.343468508671999
$ws.Range("S22").Value = [double]"7.388469271160881E-06"
$ws.Range("T22").Value = [double]"1.199631752363486E-05"

$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 0.6855039999999999
$ws.Range("H23").Value = 2.056512
$ws.Range("I23").Value = 0.002682193743693282
$ws.Range("J23").Value = 0.00341772319330189
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 0.6855039999999999
$ws.Range("N23").Value = 2.056512
$ws.Range("O23").Value = 0.002682193743693282
$ws.Range("P23").Value = 0.00341772319330189
$ws.Range("Q23").Value = 0.4699157340159998
$ws.Range("R23").Value = 4.229241606143999
$ws.Range("S23").Value = [double]"7.194163278707384E-06"
$ws.Range("T23").Value = [double]"1.168083182603367E-05"

$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 0.6855039999999999
$ws.Range("H24").Value = 2.056512
$ws.Range("I24").Value = 0.002682193743693282
$ws.Range("J24").Value = 0.00341772319330189
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 88.26666266666666
$ws.Range("N24").Value = 264.799988
$ws.Range("O24").Value = 0.3453638350486923
$ws.Range("P24").Value = 0.4400718598158738
$ws.Range("Q24").Value = 60.50715032465065
$ws.Range("R24").Value = 544.5643529218559
$ws.Range("S24").Value = 0.0009263327176655211
$ws.Range("T24").Value = 0.00150404380201221

$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 0.6855039999999999
$ws.Range("H25").Value = 2.056512
$ws.Range("I25").Value = 0.002682193743693282
$ws.Range("J25").Value = 0.00341772319330189
$ws.Range("K25").Value = 2
$ws.Range("M25").Value = 43.794673
$ws.Range("N25").Value = 87.58934600000001
$ws.Range("O25").Value = 0.1713568380748954
$ws.Range("P25").Value = 0.1455649854269483
$ws.Range("Q25").Value = 30.021423520192
$ws.Range("R25").Value = 180.128541121152
$ws.Range("S25").Value = 0.0004596122390235474
$ws.Range("T25").Value = 0.0004975008268263328

$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 88.26666266666666
$ws.Range("H26").Value = 264.799988
$ws.Range("I26").Value = 0.3453638350486923
$ws.Range("J26").Value = 0.4400718598158738
$ws.Range("K26").Value = 2
$ws.Range("M26").Value = 121.2130545
$ws.Range("N26").Value = 242.426109
$ws.Range("O26").Value = 0.4742742514031324
$ws.Range("P26").Value = 0.4028886461111009
$ws.Range("Q26").Value = 10699.07179234778
$ws.Range("R26").Value = 64194.43075408669
$ws.Range("S26").Value = 0.1637971743294334
$ws.Range("T26").Value = 0.1772999557928116

$ws.Range("E27").Value = 3
$ws.Range("G27").Value = 88.26666266666666
$ws.Range("H27").Value = 264.799988
$ws.Range("I27").Value = 0.3453638350486923
$ws.Range("J27").Value = 0.4400718598158738
$ws.Range("K27").Value = 3
$ws.Range("M27").Value = 0.9119573333333332
$ws.Range("N27").Value = 2.735872
$ws.Range("O27").Value = 0.003568245048871889
$ws.Range("P27").Value = 0.004546753526507615
$ws.Range("Q27").Value = 80.49543030772621
$ws.Range("R27").Value = 724.4588727695359
$ws.Range("S27").Value = 0.001232342794471904
$ws.Range("T27").Value = 0.002000898280534589

$ws.Range("E28").Value = 3
$ws.Range("G28").Value = 88.26666266666666
$ws.Range("H28").Value = 264.799988
$ws.Range("I28").Value = 0.3453638350486923
$ws.Range("J28").Value = 0.4400718598158738
$ws.Range("K28").Value = 3
$ws.Range("M28").Value = 0.7040186666666667
$ws.Range("N28").Value = 2.112056
$ws.Range("O28").Value = 0.002754636680714656
$ws.Range("P28").Value = 0.003510031926267592
$ws.Range("Q28").Value = 62.14137816170311
$ws.Range("R28").Value = 559.272403455328
$ws.Range("S28").Value = 0.0009513518882174138
$ws.Range("T28").Value = 0.001544666277805673

$ws.Range("E29").Value = 3
$ws.Range("G29").Value = 88.26666266666666
$ws.Range("H29").Value = 264.799988
$ws.Range("I29").Value = 0.3453638350486923
$ws.Range("J29").Value = 0.4400718598158738
$ws.Range("K29").Value = 3
$ws.Range("M29").Value = 0.6855039999999999
$ws.Range("N29").Value = 2.056512
$ws.Range("O29").Value = 0.002682193743693282
$ws.Range("P29").Value = 0.00341772319330189
$ws.Range("Q29").Value = 60.50715032465065
$ws.Range("R29").Value = 544.5643529218559
$ws.Range("S29").Value = 0.0009263327176655211
$ws.Range("T29").Value = 0.00150404380201221

$ws.Range("E30").Value = 3
$ws.Range("G30").Value = 88.26666266666666
$ws.Range("H30").Value = 264.799988
$ws.Range("I30").Value = 0.3453638350486923
$ws.Range("J30").Value = 0.4400718598158738
$ws.Range("K30").Value = 3
$ws.Range("M30").Value = 88.26666266666666
$ws.Range("N30").Value = 264.799988
$ws.Range("O30").Value = 0.3453638350486923
$ws.Range("P30").Value = 0.4400718598158738
$ws.Range("Q30").Value = 7791.003738311127
$ws.Range("R30").Value = 70119.03364480013
$ws.Range("S30").Value = 0.1192761785595403
$ws.Range("T30").Value = 0.1936632418018021

$ws.Range("E31").Value = 3
$ws.Range("G31").Value = 88.26666266666666
$ws.Range("H31").Value = 264.799988
$ws.Range("I31").Value = 0.3453638350486923
$ws.Range("J31").Value = 0.4400718598158738
$ws.Range("K31").Value = 2
$ws.Range("M31").Value = 43.794673
$ws.Range("N31").Value = 87.58934600000001
$ws.Range("O31").Value = 0.1713568380748954
$ws.Range("P31").Value = 0.1455649854269483
$ws.Range("Q31").Value = 3865.609628287975
$ws.Range("R31").Value = 23193.65776972785
$ws.Range("S31").Value = 0.05918045475936366
$ws.Range("T31").Value = 0.0640590538609077

$ws.Range("E32").Value = 2
$ws.Range("G32").Value = 43.794673
$ws.Range("H32").Value = 87.58934600000001
$ws.Range("I32").Value = 0.1713568380748954
$ws.Range("J32").Value = 0.1455649854269483
$ws.Range("K32").Value = 2
$ws.Range("M32").Value = 121.2130545
$ws.Range("N32").Value = 242.426109
$ws.Range("O32").Value = 0.4742742514031324
$ws.Range("P32").Value = 0.4028886461111009
$ws.Range("Q32").Value = 5308.486085158679
$ws.Range("R32").Value = 21233.94434063472
$ws.Range("S32").Value = 0.08127013610077881
$ws.Range("T32").Value = 0.05864647989984531

$ws.Range("E33").Value = 2
$ws.Range("G33").Value = 43.794673
$ws.Range("H33").Value = 87.58934600000001
$ws.Range("I33").Value = 0.1713568380748954
$ws.Range("J33").Value = 0.1455649854269483
$ws.Range("K33").Value = 3
$ws.Range("M33").Value = 0.9119573333333332
$ws.Range("N33").Value = 2.735872
$ws.Range("O33").Value = 0.003568245048871889
$ws.Range("P33").Value = 0.004546753526507615
$ws.Range("Q33").Value = 39.93887320328533
$ws.Range("R33").Value = 239.633239219712
$ws.Range("S33").Value = 0.0006114431890510878
$ws.Range("T33").Value = 0.0006618481108260066

$ws.Range("E34").Value = 2
$ws.Range("G34").Value = 43.794673
$ws.Range("H34").Value = 87.58934600000001
$ws.Range("I34").Value = 0.1713568380748954
$ws.Range("J34").Value = 0.1455649854269483
$ws.Range("K34").Value = 3
$ws.Range("M34").Value = 0.7040186666666667
$ws.Range("N34").Value = 2.112056
$ws.Range("O34").Value = 0.002754636680714656
$ws.Range("P34").Value = 0.003510031926267592
$ws.Range("Q34").Value = 30.83226729256267
$ws.Range("R34").Value = 184.993603755376
$ws.Range("S34").Value = 0.0004720258316523888
$ws.Range("T34").Value = 0.0005109377461952651

$ws.Range("E35").Value = 2
$ws.Range("G35").Value = 43.794673
$ws.Range("H35").Value = 87.58934600000001
$ws.Range("I35").Value = 0.1713568380748954
$ws.Range("J35").Value = 0.1455649854269483
$ws.Range("K35").Value = 3
$ws.Range("M35").Value = 0.6855039999999999
$ws.Range("N35").Value = 2.056512
$ws.Range("O35").Value = 0.002682193743693282
$ws.Range("P35").Value = 0.00341772319330189
$ws.Range("Q35").Value = 30.021423520192
$ws.Range("R35").Value = 180.128541121152
$ws.Range("S35").Value = 0.0004596122390235474
$ws.Range("T35").Value = 0.0004975008268263328

$ws.Range("E36").Value = 2
$ws.Range("G36").Value = 43.794673
$ws.Range("H36").Value = 87.58934600000001
$ws.Range("I36").Value = 0.1713568380748954
$ws.Range("J36").Value = 0.1455649854269483
$ws.Range("K36").Value = 3
$ws.Range("M36").Value = 88.26666266666666
$ws.Range("N36").Value = 264.799988
$ws.Range("O36").Value = 0.3453638350486923
$ws.Range("P36").Value = 0.4400718598158738
$ws.Range("Q36").Value = 3865.609628287975
$ws.Range("R36").Value = 23193.65776972785
$ws.Range("S36").Value = 0.05918045475936366
$ws.Range("T36").Value = 0.0640590538609077

$ws.Range("E37").Value = 2
$ws.Range("G37").Value = 43.794673
$ws.Range("H37").Value = 87.58934600000001
$ws.Range("I37").Value = 0.1713568380748954
$ws.Range("J37").Value = 0.1455649854269483
$ws.Range("K37").Value = 2
$ws.Range("M37").Value = 43.794673
$ws.Range("N37").Value = 87.58934600000001
$ws.Range("O37").Value = 0.1713568380748954
$ws.Range("P37").Value = 0.1455649854269483
$ws.Range("Q37").Value = 1917.973383176929
$ws.Range("R37").Value = 7671.893532707717
$ws.Range("S37").Value = 0.02936316595502593
$ws.Range("T37").Value = 0.02118916498234766
